# Swap the values of row 3 and row 4 for the columns that differ between
# the two records (A, B, D, E, F, G, H, M, Q, R, S). Columns C, I, K, P, T,
# U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY are identical between
# the two rows and therefore do not need to be touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "S")

foreach ($col in $cols) {
    $addr3 = "$col`3"
    $addr4 = "$col`4"
    $val3 = $ws.Range($addr3).Value2
    $val4 = $ws.Range($addr4).Value2
    $ws.Range($addr3).Value = $val4
    $ws.Range($addr4).Value = $val3
}

# Column M only has a value on row 4 before the edit; after the edit it
# only has a value on row 3. Move it explicitly instead of swapping.
$valM4 = $ws.Range("M4").Value2
$ws.Range("M3").Value = $valM4
$ws.Range("M4").Value = $null
